$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 256, pushing existing rows 256-355 down to 257-356.
$ws.Rows.Item(256).Insert()

# Populate the newly inserted row 256 with the new record's data.
$ws.Cells.Item(256, 1).Value = 9
$ws.Cells.Item(256, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(256, 3).Value = "Metropolitana"
$ws.Cells.Item(256, 4).Value = 44900
$ws.Cells.Item(256, 5).Value = 13
$ws.Cells.Item(256, 6).Value = 100112021
$ws.Cells.Item(256, 7).Value = "Ají"
$ws.Cells.Item(256, 8).Value = "Inferno"
$ws.Cells.Item(256, 9).Value = "Primera"
$ws.Cells.Item(256, 10).Value = 70
$ws.Cells.Item(256, 11).Value = 12000
$ws.Cells.Item(256, 12).Value = 15000
$ws.Cells.Item(256, 13).Value = 13500
$ws.Cells.Item(256, 14).Value = "$/caja 10 kilos"
$ws.Cells.Item(256, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(256, 16).Value = 1350
$ws.Cells.Item(256, 17).Value = 10
$ws.Cells.Item(256, 18).Value = "Hortaliza"
